# minor bom issues corrected
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (item 5, "PC/104 Headers", H1/H2): manufacturer part number corrected
$ws.Range("F6").Value = "ESQ-126-12-G-D"

# Row 12 (item 12, "10k" resistors): manufacturer part number corrected
$ws.Range("F13").Value = "RMCF0603FT10K0"

# Update the saved cursor/selection position to match
$ws.Range("F26").Select()
